$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (old D:F -> becomes E:G), making room for the
# new "VAR(OK)" column while preserving the existing MSPE/Nugget data+styles.
$ws.Range("D1").EntireColumn.Insert()

# Insert a second new column before H (old G, i.e. the data that used to be
# in F, stays put; a brand new empty column H is created) making room for
# the new "VAR(DATA)" column.
$ws.Range("H1").EntireColumn.Insert()

# ---- Header row ----
$ws.Range("D1").Value = "VAR(OK)"
# E1 already holds "MSPE" (shifted from former D1) - leave as is.
$ws.Range("F1").Value = "S_nugget"
$ws.Range("G1").Value = "VAR(TOTAL)"
$ws.Range("H1").Value = "VAR(DATA)"

# ---- Data rows ----
# New VAR(OK) column D
$ws.Range("D2").Value = 0.3634610188657926
$ws.Range("D3").Value = 0.7088661882899594
$ws.Range("D4").Value = 0.9636897536187901
$ws.Range("D5").Value = 1.116632997595095
$ws.Range("D6").Value = 1.175168753280896
$ws.Range("D7").Value = 1.200872345810631
$ws.Range("D8").Value = 1.233435112180008
$ws.Range("D9").Value = 1.255965408022047
$ws.Range("D10").Value = 1.286243006253022

# E (ex-D, MSPE) and F (ex-E, S_nugget) retain their original values already.

# The old "Non-expl var of model" values now sitting in column G are no
# longer used - clear them so G becomes the new, empty VAR(TOTAL) column.
$ws.Range("G2:G10").ClearContents()

# New VAR(DATA) column H - constant value for every row.
$ws.Range("H2:H10").Value = 6.271678887145197

$wb.Save()
